# Update "MES 01" sheet: row 54 (period 45213-45214) gets a ride recorded
# where previously there was none (VALOR 0, BAIRRO blank) -> VALOR 10,
# BAIRRO "RAIZAL". The TOTAL formula in B69 recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")
$ws.Activate()

$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "RAIZAL"

# Leave the view scrolled/selected the way the author left it.
$ws.Range("A33").Select() | Out-Null
$ws.Range("B54").Select() | Out-Null
